# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G (header "K" in G1) holds recalculated "K" values for each
# saved trade row (rows 2-16). Update them in place to the newly
# computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 1
    4  = 1
    5  = 2
    6  = 0
    7  = 1
    8  = 1
    9  = 0
    10 = 2
    11 = 0
    12 = 1
    13 = 0
    15 = 1
    16 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
